$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The yoy_schools / yoy_authorities / yoy_users columns (F, G, H) for
# rows 14-20 were mistakenly stored as percentage points (e.g. 2.94)
# instead of fractional values (0.0294). Fix by writing the corrected
# (value / 100) fractions directly, preserving full double precision.

$ws.Range("F14").Value2 = 0.02940153096729303
$ws.Range("G14").Value2 = 0.05961754780652417
$ws.Range("H14").Value2 = 0.2328061250163025

$ws.Range("F15").Value2 = 0.03658536585365857
$ws.Range("G15").Value2 = 0.06764374295377684
$ws.Range("H15").Value2 = 0.2485887932178075

$ws.Range("F16").Value2 = 0.04054289194362282
$ws.Range("G16").Value2 = 0.06877113866967299
$ws.Range("H16").Value2 = 0.2507756835683654

$ws.Range("F17").Value2 = 0.0586376404494382
$ws.Range("G17").Value2 = 0.03205128205128216
$ws.Range("H17").Value2 = 0.1887096770378025

$ws.Range("F18").Value2 = 0.06092436974789917
$ws.Range("G18").Value2 = 0.03311965811965822
$ws.Range("H18").Value2 = 0.1841667687390272

$ws.Range("F19").Value2 = 0.06339254615116685
$ws.Range("G19").Value2 = 0.03201707577374591
$ws.Range("H19").Value2 = 0.2462859203576528

$ws.Range("F20").Value2 = 0.06184142338918641
$ws.Range("G20").Value2 = 0.03503184713375807
$ws.Range("H20").Value2 = 0.2558277891171774
